# Auto-generated edit script: updates crafting-profit calculation cells
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# reflecting refreshed market-price data from the scheduled data-fetch job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 7334.75  # H18
$ws.Cells.Item(18, 9).Value = 8579.666999999999  # I18
$ws.Cells.Item(18, 11).Value = 8579.666999999999  # K18
$ws.Cells.Item(18, 13).Value = -8295.666999999999  # M18
$ws.Cells.Item(20, 8).Value = 8761  # H20
$ws.Cells.Item(20, 10).Value = 20393.5  # J20
$ws.Cells.Item(20, 12).Value = 20393.5  # L20
$ws.Cells.Item(20, 14).Value = -20853.5  # N20
$ws.Cells.Item(29, 8).Value = 599.75  # H29
$ws.Cells.Item(29, 10).Value = 599.75  # J29
$ws.Cells.Item(29, 12).Value = 1799.25  # L29
$ws.Cells.Item(29, 14).Value = -2361.25  # N29
$ws.Cells.Item(33, 8).Value = 929.9091  # H33
$ws.Cells.Item(33, 9).Value = 581  # I33
$ws.Cells.Item(33, 11).Value = 581  # K33
$ws.Cells.Item(33, 13).Value = -352  # M33
$ws.Cells.Item(35, 8).Value = 8761  # H35
$ws.Cells.Item(35, 10).Value = 20393.5  # J35
$ws.Cells.Item(35, 12).Value = 20393.5  # L35
$ws.Cells.Item(35, 14).Value = -21151.5  # N35
$ws.Cells.Item(38, 8).Value = 527.7  # H38
$ws.Cells.Item(38, 9).Value = 527.7  # I38
$ws.Cells.Item(38, 11).Value = 1583.1  # K38
$ws.Cells.Item(38, 13).Value = -1211.1  # M38
$ws.Cells.Item(40, 8).Value = 8425  # H40
$ws.Cells.Item(40, 10).Value = 8700  # J40
$ws.Cells.Item(40, 12).Value = 8700  # L40
$ws.Cells.Item(40, 14).Value = -9050  # N40
$ws.Cells.Item(41, 8).Value = 1160  # H41
$ws.Cells.Item(41, 9).Value = 666.6667  # I41
$ws.Cells.Item(41, 10).Value = 1900  # J41
$ws.Cells.Item(41, 11).Value = 666.6667  # K41
$ws.Cells.Item(41, 12).Value = 1900  # L41
$ws.Cells.Item(41, 13).Value = -226.6667  # M41
$ws.Cells.Item(41, 14).Value = -2780  # N41
$ws.Cells.Item(46, 8).Value = 115986.445  # H46
$ws.Cells.Item(46, 9).Value = 1450  # I46
$ws.Cells.Item(46, 10).Value = 148711.14  # J46
$ws.Cells.Item(46, 11).Value = 4350  # K46
$ws.Cells.Item(46, 12).Value = 446133.42  # L46
$ws.Cells.Item(46, 13).Value = -4231  # M46
$ws.Cells.Item(46, 14).Value = -446371.42  # N46
$ws.Cells.Item(48, 8).Value = 5026.1333  # H48
$ws.Cells.Item(48, 10).Value = 5116.4316  # J48
$ws.Cells.Item(48, 12).Value = 15349.2948  # L48
$ws.Cells.Item(48, 14).Value = -15933.2948  # N48
$ws.Cells.Item(49, 8).Value = 4348.7144  # H49
$ws.Cells.Item(49, 9).Value = 17  # I49
$ws.Cells.Item(49, 10).Value = 5070.6665  # J49
$ws.Cells.Item(49, 11).Value = 51  # K49
$ws.Cells.Item(49, 12).Value = 15211.9995  # L49
$ws.Cells.Item(49, 14).Value = -15483.9995  # N49
$ws.Cells.Item(49, 13).Value = 85  # M49
$ws.Cells.Item(56, 8).Value = 5026.1333  # H56
$ws.Cells.Item(56, 10).Value = 5116.4316  # J56
$ws.Cells.Item(56, 12).Value = 15349.2948  # L56
$ws.Cells.Item(56, 14).Value = -16417.2948  # N56
$ws.Cells.Item(60, 8).Value = 115986.445  # H60
$ws.Cells.Item(60, 9).Value = 1450  # I60
$ws.Cells.Item(60, 10).Value = 148711.14  # J60
$ws.Cells.Item(60, 11).Value = 4350  # K60
$ws.Cells.Item(60, 12).Value = 446133.42  # L60
$ws.Cells.Item(60, 13).Value = -3866  # M60
$ws.Cells.Item(60, 14).Value = -447101.42  # N60
$ws.Cells.Item(107, 8).Value = 7234.143  # H107
$ws.Cells.Item(107, 9).Value = 6943.263  # I107
$ws.Cells.Item(107, 10).Value = 9997.5  # J107
$ws.Cells.Item(107, 11).Value = 6943.263  # K107
$ws.Cells.Item(107, 12).Value = 9997.5  # L107
$ws.Cells.Item(107, 13).Value = -5023.263  # M107
$ws.Cells.Item(107, 14).Value = -13837.5  # N107
$ws.Cells.Item(137, 8).Value = 7001.091  # H137
$ws.Cells.Item(137, 9).Value = 9610.286  # I137
$ws.Cells.Item(137, 10).Value = 2435  # J137
$ws.Cells.Item(137, 11).Value = 28830.858  # K137
$ws.Cells.Item(137, 12).Value = 7305  # L137
$ws.Cells.Item(137, 13).Value = -26280.858  # M137
$ws.Cells.Item(137, 14).Value = -12405  # N137
$ws.Cells.Item(138, 8).Value = 3599.0476  # H138
$ws.Cells.Item(138, 9).Value = 1558  # I138
$ws.Cells.Item(138, 10).Value = 5454.5454  # J138
$ws.Cells.Item(138, 11).Value = 4674  # K138
$ws.Cells.Item(138, 12).Value = 16363.6362  # L138
$ws.Cells.Item(138, 13).Value = 466  # M138
$ws.Cells.Item(138, 14).Value = -26643.6362  # N138
$ws.Cells.Item(141, 8).Value = 4515.838  # H141
$ws.Cells.Item(141, 9).Value = 3186.2666  # I141
$ws.Cells.Item(141, 10).Value = 10214  # J141
$ws.Cells.Item(141, 11).Value = 9558.799800000001  # K141
$ws.Cells.Item(141, 12).Value = 30642  # L141
$ws.Cells.Item(141, 13).Value = -4378.799800000001  # M141
$ws.Cells.Item(141, 14).Value = -41002  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 46018.824  # H2
$ws.Cells.Item(2, 9).Value = 2126.6365  # I2
$ws.Cells.Item(2, 11).Value = 2126.6365  # K2
$ws.Cells.Item(2, 13).Value = -2013.6365  # M2
$ws.Cells.Item(32, 8).Value = 3673.9246  # H32
$ws.Cells.Item(32, 9).Value = 3673.9246  # I32
$ws.Cells.Item(32, 11).Value = 3673.9246  # K32
$ws.Cells.Item(32, 13).Value = -3386.9246  # M32
$ws.Cells.Item(37, 8).Value = 5428.4287  # H37
$ws.Cells.Item(37, 9).Value = 5428.4287  # I37
$ws.Cells.Item(37, 11).Value = 5428.4287  # K37
$ws.Cells.Item(37, 13).Value = -5155.4287  # M37
$ws.Cells.Item(102, 8).Value = 14462.117  # H102
$ws.Cells.Item(102, 9).Value = 16418.715  # I102
$ws.Cells.Item(102, 11).Value = 16418.715  # K102
$ws.Cells.Item(102, 13).Value = -14796.715  # M102
$ws.Cells.Item(110, 8).Value = 3189.1365  # H110
$ws.Cells.Item(110, 9).Value = 2546.6428  # I110
$ws.Cells.Item(110, 11).Value = 2546.6428  # K110
$ws.Cells.Item(110, 13).Value = -501.6428000000001  # M110
$ws.Cells.Item(116, 8).Value = 46018.824  # H116
$ws.Cells.Item(116, 9).Value = 2126.6365  # I116
$ws.Cells.Item(116, 11).Value = 2126.6365  # K116
$ws.Cells.Item(116, 13).Value = 167.3634999999999  # M116
$ws.Cells.Item(132, 8).Value = 4811.1396  # H132
$ws.Cells.Item(132, 9).Value = 4762  # I132
$ws.Cells.Item(132, 11).Value = 14286  # K132
$ws.Cells.Item(132, 13).Value = -11756  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 46018.824  # H3
$ws.Cells.Item(3, 9).Value = 2126.6365  # I3
$ws.Cells.Item(3, 11).Value = 2126.6365  # K3
$ws.Cells.Item(3, 13).Value = -2012.6365  # M3
$ws.Cells.Item(10, 8).Value = 3361.25  # H10
$ws.Cells.Item(10, 9).Value = 445  # I10
$ws.Cells.Item(10, 10).Value = 4333.3335  # J10
$ws.Cells.Item(10, 11).Value = 445  # K10
$ws.Cells.Item(10, 12).Value = 4333.3335  # L10
$ws.Cells.Item(10, 13).Value = -305  # M10
$ws.Cells.Item(10, 14).Value = -4613.3335  # N10
$ws.Cells.Item(35, 8).Value = 77497.5  # H35
$ws.Cells.Item(35, 10).Value = 77497.5  # J35
$ws.Cells.Item(35, 12).Value = 77497.5  # L35
$ws.Cells.Item(35, 14).Value = -78117.5  # N35
$ws.Cells.Item(99, 8).Value = 23279.723  # H99
$ws.Cells.Item(99, 9).Value = 26298.928  # I99
$ws.Cells.Item(99, 10).Value = 12712.5  # J99
$ws.Cells.Item(99, 11).Value = 26298.928  # K99
$ws.Cells.Item(99, 12).Value = 12712.5  # L99
$ws.Cells.Item(99, 13).Value = -24800.928  # M99
$ws.Cells.Item(99, 14).Value = -15708.5  # N99
$ws.Cells.Item(103, 8).Value = 21300  # H103
$ws.Cells.Item(103, 10).Value = 21300  # J103
$ws.Cells.Item(103, 12).Value = 21300  # L103
$ws.Cells.Item(103, 14).Value = -23644  # N103
$ws.Cells.Item(107, 8).Value = 2362.45  # H107
$ws.Cells.Item(107, 9).Value = 2373.4707  # I107
$ws.Cells.Item(107, 11).Value = 2373.4707  # K107
$ws.Cells.Item(107, 13).Value = -453.4706999999999  # M107
$ws.Cells.Item(132, 8).Value = 82111.60000000001  # H132
$ws.Cells.Item(132, 10).Value = 82111.60000000001  # J132
$ws.Cells.Item(132, 12).Value = 82111.60000000001  # L132
$ws.Cells.Item(132, 14).Value = -92231.60000000001  # N132
$ws.Cells.Item(134, 8).Value = 6916.8  # H134
$ws.Cells.Item(134, 9).Value = 7142.317  # I134
$ws.Cells.Item(134, 11).Value = 21426.951  # K134
$ws.Cells.Item(134, 13).Value = -18891.951  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1656.2858  # H16
$ws.Cells.Item(16, 10).Value = 1498  # J16
$ws.Cells.Item(16, 12).Value = 1498  # L16
$ws.Cells.Item(16, 14).Value = -2072  # N16
$ws.Cells.Item(31, 8).Value = 4139.5454  # H31
$ws.Cells.Item(31, 9).Value = 3464.818  # I31
$ws.Cells.Item(31, 10).Value = 5489  # J31
$ws.Cells.Item(31, 11).Value = 3464.818  # K31
$ws.Cells.Item(31, 12).Value = 5489  # L31
$ws.Cells.Item(31, 13).Value = -3169.818  # M31
$ws.Cells.Item(31, 14).Value = -6079  # N31
$ws.Cells.Item(34, 8).Value = 4139.5454  # H34
$ws.Cells.Item(34, 9).Value = 3464.818  # I34
$ws.Cells.Item(34, 10).Value = 5489  # J34
$ws.Cells.Item(34, 11).Value = 3464.818  # K34
$ws.Cells.Item(34, 12).Value = 5489  # L34
$ws.Cells.Item(34, 13).Value = -3262.818  # M34
$ws.Cells.Item(34, 14).Value = -5893  # N34
$ws.Cells.Item(99, 8).Value = 420883.34  # H99
$ws.Cells.Item(99, 10).Value = 5775  # J99
$ws.Cells.Item(99, 12).Value = 5775  # L99
$ws.Cells.Item(99, 14).Value = -8771  # N99
$ws.Cells.Item(113, 8).Value = 1656.2858  # H113
$ws.Cells.Item(113, 10).Value = 1498  # J113
$ws.Cells.Item(113, 12).Value = 1498  # L113
$ws.Cells.Item(113, 14).Value = -5838  # N113
$ws.Cells.Item(126, 8).Value = 420883.34  # H126
$ws.Cells.Item(126, 10).Value = 5775  # J126
$ws.Cells.Item(126, 12).Value = 17325  # L126
$ws.Cells.Item(126, 14).Value = -22265  # N126
$ws.Cells.Item(132, 8).Value = 19984.773  # H132
$ws.Cells.Item(132, 9).Value = 1450.2106  # I132
$ws.Cells.Item(132, 11).Value = 4350.6318  # K132
$ws.Cells.Item(132, 13).Value = -1820.6318  # M132
$ws.Cells.Item(134, 8).Value = 2105.1887  # H134
$ws.Cells.Item(134, 9).Value = 770.4651  # I134
$ws.Cells.Item(134, 10).Value = 7844.5  # J134
$ws.Cells.Item(134, 11).Value = 2311.3953  # K134
$ws.Cells.Item(134, 12).Value = 23533.5  # L134
$ws.Cells.Item(134, 13).Value = 223.6046999999999  # M134
$ws.Cells.Item(134, 14).Value = -28603.5  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 180.1579  # H12
$ws.Cells.Item(12, 9).Value = 157  # I12
$ws.Cells.Item(12, 10).Value = 201  # J12
$ws.Cells.Item(12, 11).Value = 471  # K12
$ws.Cells.Item(12, 12).Value = 603  # L12
$ws.Cells.Item(12, 13).Value = -298  # M12
$ws.Cells.Item(12, 14).Value = -949  # N12
$ws.Cells.Item(46, 8).Value = 5047.9  # H46
$ws.Cells.Item(46, 9).Value = 0  # I46
$ws.Cells.Item(46, 10).Value = 5047.9  # J46
$ws.Cells.Item(46, 11).Value = 0  # K46
$ws.Cells.Item(46, 12).Value = 15143.7  # L46
$ws.Cells.Item(46, 13).Value = ""  # M46
$ws.Cells.Item(46, 14).Value = -15325.7  # N46
$ws.Cells.Item(109, 8).Value = 6079.8  # H109
$ws.Cells.Item(109, 9).Value = 400  # I109
$ws.Cells.Item(109, 11).Value = 1200  # K109
$ws.Cells.Item(109, 13).Value = -160  # M109

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(59, 8).Value = 11499.5  # H59
$ws.Cells.Item(59, 10).Value = 10999  # J59
$ws.Cells.Item(59, 12).Value = 10999  # L59
$ws.Cells.Item(59, 14).Value = -12165  # N59
$ws.Cells.Item(70, 8).Value = 13556.5  # H70
$ws.Cells.Item(70, 9).Value = 11761.8  # I70
$ws.Cells.Item(70, 10).Value = 14372.272  # J70
$ws.Cells.Item(70, 11).Value = 11761.8  # K70
$ws.Cells.Item(70, 12).Value = 14372.272  # L70
$ws.Cells.Item(70, 13).Value = -11491.8  # M70
$ws.Cells.Item(70, 14).Value = -14912.272  # N70
$ws.Cells.Item(73, 8).Value = 13556.5  # H73
$ws.Cells.Item(73, 9).Value = 11761.8  # I73
$ws.Cells.Item(73, 10).Value = 14372.272  # J73
$ws.Cells.Item(73, 11).Value = 11761.8  # K73
$ws.Cells.Item(73, 12).Value = 14372.272  # L73
$ws.Cells.Item(73, 13).Value = -10825.8  # M73
$ws.Cells.Item(73, 14).Value = -16244.272  # N73
$ws.Cells.Item(97, 8).Value = 6575  # H97
$ws.Cells.Item(97, 9).Value = 7607.9062  # I97
$ws.Cells.Item(97, 11).Value = 7607.9062  # K97
$ws.Cells.Item(97, 13).Value = -7111.9062  # M97
$ws.Cells.Item(107, 8).Value = 519.1852  # H107
$ws.Cells.Item(107, 9).Value = 549.1429000000001  # I107
$ws.Cells.Item(107, 10).Value = 414.33334  # J107
$ws.Cells.Item(107, 11).Value = 549.1429000000001  # K107
$ws.Cells.Item(107, 12).Value = 414.33334  # L107
$ws.Cells.Item(107, 13).Value = 1370.8571  # M107
$ws.Cells.Item(107, 14).Value = -4254.33334  # N107
$ws.Cells.Item(126, 8).Value = 22583.666  # H126
$ws.Cells.Item(126, 9).Value = 28239.625  # I126
$ws.Cells.Item(126, 11).Value = 84718.875  # K126
$ws.Cells.Item(126, 13).Value = -82248.875  # M126
$ws.Cells.Item(132, 8).Value = 2628.9678  # H132
$ws.Cells.Item(132, 9).Value = 2259.625  # I132
$ws.Cells.Item(132, 11).Value = 6778.875  # K132
$ws.Cells.Item(132, 13).Value = -4248.875  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2974.0715  # H16
$ws.Cells.Item(16, 9).Value = 2776  # I16
$ws.Cells.Item(16, 11).Value = 2776  # K16
$ws.Cells.Item(16, 13).Value = -2606  # M16
$ws.Cells.Item(46, 8).Value = 3716.5715  # H46
$ws.Cells.Item(46, 10).Value = 4866.3335  # J46
$ws.Cells.Item(46, 12).Value = 4866.3335  # L46
$ws.Cells.Item(46, 14).Value = -5242.3335  # N46
$ws.Cells.Item(61, 8).Value = 5256.4517  # H61
$ws.Cells.Item(61, 9).Value = 4056.238  # I61
$ws.Cells.Item(61, 11).Value = 4056.238  # K61
$ws.Cells.Item(61, 13).Value = -3854.238  # M61
$ws.Cells.Item(68, 8).Value = 3452.0833  # H68
$ws.Cells.Item(68, 9).Value = 2203.1428  # I68
$ws.Cells.Item(68, 10).Value = 5200.6  # J68
$ws.Cells.Item(68, 11).Value = 2203.1428  # K68
$ws.Cells.Item(68, 12).Value = 5200.6  # L68
$ws.Cells.Item(68, 13).Value = -1454.1428  # M68
$ws.Cells.Item(68, 14).Value = -6698.6  # N68
$ws.Cells.Item(71, 8).Value = 3452.0833  # H71
$ws.Cells.Item(71, 9).Value = 2203.1428  # I71
$ws.Cells.Item(71, 10).Value = 5200.6  # J71
$ws.Cells.Item(71, 11).Value = 11015.714  # K71
$ws.Cells.Item(71, 12).Value = 26003  # L71
$ws.Cells.Item(71, 13).Value = -7271.714  # M71
$ws.Cells.Item(71, 14).Value = -33491  # N71
$ws.Cells.Item(93, 8).Value = 5927.3887  # H93
$ws.Cells.Item(93, 9).Value = 7092  # I93
$ws.Cells.Item(93, 10).Value = 1851.25  # J93
$ws.Cells.Item(93, 11).Value = 7092  # K93
$ws.Cells.Item(93, 12).Value = 1851.25  # L93
$ws.Cells.Item(93, 13).Value = -5844  # M93
$ws.Cells.Item(93, 14).Value = -4347.25  # N93
$ws.Cells.Item(100, 8).Value = 3105.625  # H100
$ws.Cells.Item(100, 9).Value = 3132.5  # I100
$ws.Cells.Item(100, 11).Value = 3132.5  # K100
$ws.Cells.Item(100, 13).Value = -2591.5  # M100
$ws.Cells.Item(113, 8).Value = 5256.4517  # H113
$ws.Cells.Item(113, 9).Value = 4056.238  # I113
$ws.Cells.Item(113, 11).Value = 4056.238  # K113
$ws.Cells.Item(113, 13).Value = -1886.238  # M113
$ws.Cells.Item(122, 8).Value = 6669.4165  # H122
$ws.Cells.Item(122, 9).Value = 5503.6787  # I122
$ws.Cells.Item(122, 11).Value = 16511.0361  # K122
$ws.Cells.Item(122, 13).Value = -14061.0361  # M122
$ws.Cells.Item(132, 8).Value = 624983.5600000001  # H132
$ws.Cells.Item(132, 9).Value = 933288.4399999999  # I132
$ws.Cells.Item(132, 11).Value = 2799865.32  # K132
$ws.Cells.Item(132, 13).Value = -2797335.32  # M132
$ws.Cells.Item(136, 8).Value = 7578.478  # H136
$ws.Cells.Item(136, 9).Value = 2402  # I136
$ws.Cells.Item(136, 10).Value = 14307.9  # J136
$ws.Cells.Item(136, 11).Value = 7206  # K136
$ws.Cells.Item(136, 12).Value = 42923.7  # L136
$ws.Cells.Item(136, 13).Value = -4656  # M136
$ws.Cells.Item(136, 14).Value = -48023.7  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 39995  # H21
$ws.Cells.Item(21, 10).Value = 39995  # J21
$ws.Cells.Item(21, 12).Value = 39995  # L21
$ws.Cells.Item(21, 14).Value = -40465  # N21
$ws.Cells.Item(35, 8).Value = 39995  # H35
$ws.Cells.Item(35, 10).Value = 39995  # J35
$ws.Cells.Item(35, 12).Value = 39995  # L35
$ws.Cells.Item(35, 14).Value = -40575  # N35
$ws.Cells.Item(51, 8).Value = 33448.168  # H51
$ws.Cells.Item(51, 9).Value = 6933  # I51
$ws.Cells.Item(51, 11).Value = 6933  # K51
$ws.Cells.Item(51, 13).Value = -6423  # M51
$ws.Cells.Item(96, 8).Value = 1820.7778  # H96
$ws.Cells.Item(96, 10).Value = 1875  # J96
$ws.Cells.Item(96, 12).Value = 1875  # L96
$ws.Cells.Item(96, 14).Value = -4621  # N96
$ws.Cells.Item(107, 8).Value = 12477.704  # H107
$ws.Cells.Item(107, 9).Value = 1501.1904  # I107
$ws.Cells.Item(107, 10).Value = 50895.5  # J107
$ws.Cells.Item(107, 11).Value = 4503.5712  # K107
$ws.Cells.Item(107, 12).Value = 152686.5  # L107
$ws.Cells.Item(107, 13).Value = -2583.5712  # M107
$ws.Cells.Item(107, 14).Value = -156526.5  # N107
$ws.Cells.Item(113, 8).Value = 1029  # H113
$ws.Cells.Item(113, 9).Value = 499.7353  # I113
$ws.Cells.Item(113, 10).Value = 10026.5  # J113
$ws.Cells.Item(113, 11).Value = 1499.2059  # K113
$ws.Cells.Item(113, 12).Value = 30079.5  # L113
$ws.Cells.Item(113, 13).Value = 670.7941000000001  # M113
$ws.Cells.Item(113, 14).Value = -34419.5  # N113
$ws.Cells.Item(122, 8).Value = 3790.1143  # H122
$ws.Cells.Item(122, 9).Value = 1779.7273  # I122
$ws.Cells.Item(122, 11).Value = 5339.1819  # K122
$ws.Cells.Item(122, 13).Value = -2889.1819  # M122
$ws.Cells.Item(132, 8).Value = 9824  # H132
$ws.Cells.Item(132, 9).Value = 10131.548  # I132
$ws.Cells.Item(132, 11).Value = 30394.644  # K132
$ws.Cells.Item(132, 13).Value = -27864.644  # M132
$ws.Cells.Item(136, 8).Value = 455217.25  # H136
$ws.Cells.Item(136, 9).Value = 593284.0600000001  # I136
$ws.Cells.Item(136, 11).Value = 1779852.18  # K136
$ws.Cells.Item(136, 13).Value = -1777302.18  # M136

